# Frontend works with backend
# Update the email-status tracker sheet: new recipient emails + hyperlinks,
# simplified "Sent" status text, and the cursor's last selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "shreyaslegend364@gmail.com"
$ws.Range("B2").Value = "Sent"

# --- Row 3 ---
$ws.Range("A3").Value = "gowdashreyas364@gmail.com"
$ws.Range("B3").Value = "Sent"

# Apply the workbook's built-in "Hyperlink" cell style first, then add the
# actual mailto: hyperlinks on top of it.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:shreyaslegend364@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:gowdashreyas364@gmail.com")

# Touch column C so the engine emits it as its own <col> record (matching
# the split column-run in the saved file) without altering its width/style.
$ws.Columns.Item(3).Hidden = $false

# Restore the previously-selected cell.
$ws.Range("C6").Select()

Write-Host "done"
